$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.621.29"
$ws.Range("E2").Value = "  +1.71%  "

$ws.Range("D3").Value = "1.599.32"
$ws.Range("E3").Value = "  +1.35%  "

$ws.Range("E4").Value = "  +0.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("E7").Value = "  +0.50%  "

$ws.Range("E8").Value = "  +3.93%  "

$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("E10").Value = "  +1.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.04%  "

$ws.Range("D12").Value = "1.827.63"
$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("D13").Value = "1.601.17"
$ws.Range("E13").Value = "  +0.90%  "

$ws.Range("D14").Value = "29.620.01"
$ws.Range("E14").Value = "  +1.68%  "

$ws.Range("E15").Value = "  +3.12%  "

$ws.Range("E16").Value = "  +0.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.95"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241.50"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.69%  "

$ws.Range("E19").Value = "  +2.49%  "

$ws.Range("D20").Value = "0.0₃0693"
$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("E22").Value = "  -0.30%  "

$ws.Range("E23").Value = "  +0.55%  "

$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.45"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.35"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.86%  "

$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.39"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.28%  "

$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("E30").Value = "  +2.80%  "

$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("E33").Value = "  +3.19%  "

$ws.Range("D34").Value = "1.425.44"
$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("E35").Value = "  +2.12%  "

$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.88"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.07%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.03"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.30"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.49%  "

$ws.Range("E39").Value = "  +2.06%  "

$ws.Range("E40").Value = "  +3.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.96"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.03%  "

$ws.Range("E42").Value = "  +4.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "54.48"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.46%  "

$ws.Range("E44").Value = "  +2.14%  "

$ws.Range("E45").Value = "  +0.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.993"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +17.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.78"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.92%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.31"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.46%  "

$ws.Range("D49").Value = "1.738.34"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.04"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("D51").Value = "0.0₆0106"
$ws.Range("E51").Value = "  +7.56%  "
